# Update Betfair Back/Lay odds data for 2026-01-23 per the day's refreshed snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.76
$ws.Range("G2").Value = 1.79
$ws.Range("H2").Value = 4.4
$ws.Range("I2").Value = 4.6
$ws.Range("J2").Value = 4.6
$ws.Range("K2").Value = 4.8
$ws.Range("L2").Value = 1.21
$ws.Range("N2").Value = 7.6
$ws.Range("P2").Value = 3.15
$ws.Range("R2").Value = 1.85
$ws.Range("S2").Value = 1.86
$ws.Range("T2").Value = 1.45
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 1.28
$ws.Range("W2").Value = 2.26
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 980
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 980
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 980
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.25
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 1.12
$ws.Range("S3").Value = 1.01
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 65
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.36
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1.36
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1.02
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 1.25
$ws.Range("P4").Value = 1.25
$ws.Range("Q4").Value = 1.39
$ws.Range("R4").Value = 1.16
$ws.Range("S4").Value = 1.01
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.84
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("Y4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AN4").Value = 1000

# Row 5
$ws.Range("H5").Value = 1.33
$ws.Range("J5").Value = 5.9
$ws.Range("N5").Value = 4.8
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 2.44
$ws.Range("Q5").Value = 1.49
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = 2.26
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.98
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 6
$ws.Range("F6").Value = 2.72
$ws.Range("G6").Value = 2.74
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 3.1
$ws.Range("P6").Value = 1.68
$ws.Range("Q6").Value = 2.36
$ws.Range("T6").Value = 1.98
$ws.Range("U6").Value = 1.95
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 11
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 10.5
$ws.Range("AC6").Value = 7.4
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 15.5
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 980
$ws.Range("AO6").Value = 980

# Row 7
$ws.Range("H7").Value = 19
$ws.Range("P7").Value = 2.52
$ws.Range("Q7").Value = 1.58
$ws.Range("S7").Value = 2.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 1.57
$ws.Range("X7").Value = 980
$ws.Range("Y7").Value = 350
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AF7").Value = 7.4
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 360
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 570
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 3.75

# Row 8
$ws.Range("F8").Value = 2.62
$ws.Range("G8").Value = 2.7
$ws.Range("H8").Value = 2.94
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 3.45
$ws.Range("K8").Value = 3.5
$ws.Range("M8").Value = 1.07
$ws.Range("Q8").Value = 2.04
$ws.Range("X8").Value = 14
$ws.Range("Y8").Value = 12.5
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 11.5
$ws.Range("AE8").Value = 980
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 14
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 980
$ws.Range("AO8").Value = 980

# Row 9
$ws.Range("F9").Value = 2.48
$ws.Range("G9").Value = 2.76
$ws.Range("H9").Value = 3.4
$ws.Range("K9").Value = 3.2

# Row 10
$ws.Range("F10").Value = 2.24
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 1000
$ws.Range("K10").Value = 1000
